$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 49 : Minifit jr R/A 4 pin connector
# ---------------------------------------------------------------------------
$ws.Range("E49").Value2 = "Connector"
$ws.Range("D49").Value2 = "WM13195-ND"
$ws.Range("B49").Value2 = "https://www.digikey.com/en/products/detail/molex/0353180420/3185063?s=N4IgTCBcDaIOoFkCMBmJBOArAWgHIBEQBdAXyA"
$ws.Range("G49").Value2 = "Minifit jr R/A 4 pin"

# ---------------------------------------------------------------------------
# Row 50 : Minifit jr R/A 2 pin connector
# ---------------------------------------------------------------------------
$ws.Range("B50").Value2 = "https://www.digikey.com/en/products/detail/molex/0039300020/930320"
$ws.Range("D50").Value2 = "WM21351-ND"
$ws.Range("E50").Value2 = "Connector"
$ws.Range("G50").Value2 = "Minifit jr R/A 2 pin"

# ---------------------------------------------------------------------------
# Notes column (K) - added after the rest of the two rows
# ---------------------------------------------------------------------------
$ws.Range("K50").Value2 = "Weren't listed in the inventory report but we may have some"
$ws.Range("K49").Value2 = "This and the 2 pin below also need their respective connectors"

# ---------------------------------------------------------------------------
# Distributor column reuses the existing "Digikey" shared string
# ---------------------------------------------------------------------------
$ws.Range("C49").Value2 = "Digikey"
$ws.Range("C50").Value2 = "Digikey"

# ---------------------------------------------------------------------------
# Manufacturer part numbers (numeric, left aligned)
# ---------------------------------------------------------------------------
$ws.Range("F49").Value2 = 353180420
$ws.Range("F49").Style = "Normal"
$ws.Range("F49").HorizontalAlignment = -4131

$ws.Range("F50").Value2 = 39300020
$ws.Range("F50").HorizontalAlignment = -4131

# ---------------------------------------------------------------------------
# Quantities
# ---------------------------------------------------------------------------
$ws.Range("H49").Value2 = 2
$ws.Range("I49").Value2 = 5

$ws.Range("H50").Value2 = 2
$ws.Range("I50").Value2 = 0

# ---------------------------------------------------------------------------
# Difference formulas (J column, shared formula group)
# ---------------------------------------------------------------------------
$ws.Range("J49").Formula = "=I49-H49"
$ws.Range("J50").Formula = "=I50-H50"

# ---------------------------------------------------------------------------
# Fix up styles for D49 (no explicit border/format, unlike D50)
# ---------------------------------------------------------------------------
$ws.Range("D49").Style = "Normal"

# ---------------------------------------------------------------------------
# Copy matching cell formats from the row above (B/C/D/E/G/H/I/J) so that the
# new rows visually match the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("B48").Copy()
$ws.Range("B49:B50").PasteSpecial(-4122)

$ws.Range("C48").Copy()
$ws.Range("C49:C50").PasteSpecial(-4122)

$ws.Range("E48").Copy()
$ws.Range("E49:E50").PasteSpecial(-4122)

$ws.Range("G48").Copy()
$ws.Range("G49:G50").PasteSpecial(-4122)

$ws.Range("H48:J48").Copy()
$ws.Range("H49:J50").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Hyperlinks for the new distributor links
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B49"), "https://www.digikey.com/en/products/detail/molex/0353180420/3185063?s=N4IgTCBcDaIOoFkCMBmJBOArAWgHIBEQBdAXyA")
$ws.Hyperlinks.Add($ws.Range("B50"), "https://www.digikey.com/en/products/detail/molex/0039300020/930320")

$ws.Range("B48").Copy()
$ws.Range("B49:B50").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Update the view state to reflect the new scroll/selection position
# ---------------------------------------------------------------------------
$ws.Range("B51").Select()

$wb.Save()
